$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new data points to row 2 (Job_Id = JD_001)
$ws.Range("F2").Value = "Created"
$ws.Range("G2").Value = "Yes"

# Match the saved selection state from the authored file
$ws.Range("H7").Select()
